$d = $word.ActiveDocument
$d.Content.Find.Execute("Project #1", $true, $false, $false, $false, $false, $true, 1, $false, "Project #2", 2)
